# Apply updated "dSF" (column F) values for this data repull / mean calculation pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new dSF (column F) value
$updates = @{
    2  = 1
    3  = 3
    7  = -1
    8  = 3
    10 = -3
    13 = -2
    20 = -7
    21 = 3
    22 = -2
    23 = 2
    25 = -3
    27 = 0
    28 = -4
    29 = -3
    30 = 2
    31 = -16
    32 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
